$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")

# Insert a new column before column E (shifts E:N -> F:O), matching the
# "add Expected Meganav Breadcrumb column" edit.
$ws.Columns("E:E").Insert()

# Populate the new column. Set E2 before E1 so the shared-string table
# gains "Home/Tableware/Serveware & Flatware" before
# "Expected Meganav Breadcrumb" (matches authoring order in the diff).
$ws.Range("E2").Value = "Home/Tableware/Serveware & Flatware"
$ws.Range("E1").Value = "Expected Meganav Breadcrumb"

# New column mirrors the wide "breadcrumb" style column width used
# elsewhere in the sheet (~62.57 chars; 61.71 is the closest input that
# rounds to that stored width).
$ws.Columns("E:E").ColumnWidth = 61.71

# The hyperlink that lived on (old) N2 needs to end up on (new) O2 - the
# plain column insert doesn't retarget it, so rebuild it in place.
$url = "https://surlatable.testrail.net/index.php?/cases/view/12080&group_by=cases:section_id&group_order=asc&display_deleted_cases=0&group_id=1961"
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("O2"), $url)
$ws.Range("O2").Style = "Hyperlink"

# Update the active selection to the new column (matches the saved view).
$ws.Range("E1").Select()
